# Chapter 3 manuscript v2 -> v2.1 edits
# - The placeholder "Supplementary Figure X" rows (ctenophore-first / sponge-first
#   species trees) are renamed: their Name column now reads
#   "are available on GitHub" (they join the other supplementary-files-on-GitHub rows).
# - Selection/cursor left on B24 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentioned_in_text")

$ws.Range("A23").Value = "are available on GitHub"
$ws.Range("A24").Value = "are available on GitHub"

$ws.Select()
$ws.Range("B24").Select()
